# Unify the conception of DataNode, DataTable, Entity.
# Substantive content changes from the commit:
#  - Sheet renamed from "Property1" to "DataNode"
#  - Active selection moved to E23 (bottom pane of the frozen view)
#  - Minor column width touch-ups on columns A and E

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet.
$ws.Name = "DataNode"

# Nudge the stored column widths for columns A and E closer to their
# new values (12.625 / 13.875 characters).
$ws.Columns.Item(1).ColumnWidth = 11.857142857142858
$ws.Columns.Item(5).ColumnWidth = 13.142857142857142

# Move/restore the active selection to E23 in the frozen bottom-left pane.
$ws.Range("E23").Select()
